$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Hide helper columns F and G ---
$ws.Columns("F").Hidden = $true
$ws.Columns("G").Hidden = $true

# --- Update product descriptions in column H (rows 2-20) and wrap text ---
# Row 2 keeps its existing description text; only formatting (wrap) changes.
$ws.Range("H2").WrapText = $true

$ws.Range("H3").Value = 'Enjoy the delicious taste of Sphinx Bagels with melted butter or cream cheese for a blissful breakfast or anytime snack.'
$ws.Range("H3").WrapText = $true
$ws.Rows(3).RowHeight = 32

$ws.Range("H4").Value = 'Great Value Original EnglishSweet, tangy or savory, and sometimes a combination—with so much variety, there''s plenty to explore'
$ws.Range("H4").WrapText = $true
$ws.Rows(4).RowHeight = 32

$ws.Range("H5").Value = 'The extra thickness Great Cranberry signature muffins will stand up to any topping, and the tiny air pockets make for a crisp bite when they come out of the toaster.'
$ws.Range("H5").WrapText = $true
$ws.Rows(5).RowHeight = 32

$ws.Range("H6").Value = 'Our gourmet Blueberry muffins are freezer friendly, so it’s a good idea to stock up before the brunch season is in full swing'
$ws.Range("H6").WrapText = $true
$ws.Rows(6).RowHeight = 32

$ws.Range("H7").Value = 'For a delicious and satisfying breakfast or brunch treat or snack, top them with one of our gourmet fruit preserves, indulgent butters, and other fabulous toppings.'
$ws.Range("H7").WrapText = $true
$ws.Rows(7).RowHeight = 32

$ws.Range("H8").Value = 'Pumpernickel Bread:Also delicious when toasted High fiber,Cholesterol free,Natural ingredients,Long shelf life,Wheat free ,Kosher'
$ws.Range("H8").WrapText = $true
$ws.Rows(8).RowHeight = 32

$ws.Range("H9").Value = 'Great White Bread  : Gluten Free Wheat Free Corn Free Dairy Free Soy Free'
$ws.Range("H9").WrapText = $true

$ws.Range("H10").Value = 'Very good with cold cuts, fish, poultry, soups. Makes tasty sandwiches. Suitable for sweet topping & Sour Taste.  INGREDIENTS: Dark rye flour, Austrian sourdough starter, filtered water, sea salt.'
$ws.Range("H10").WrapText = $true
$ws.Rows(10).RowHeight = 48

$ws.Range("H11").Value = 'It''s a perfect combination of whether you are interested in finding your flavor or you are just interested in mixing it up'
$ws.Range("H11").WrapText = $true
$ws.Rows(11).RowHeight = 32

$ws.Range("H12").Value = 'Flavors in variety pack include diet black cherry, diet cherry vanilla cream, diet root beer and diet tangerine lime'
$ws.Range("H12").WrapText = $true
$ws.Rows(12).RowHeight = 32

$ws.Range("H13").Value = 'No sugar added Excellent source of Vitamin C for daily value. Provides one cup of fruit, so kids get the fruit they need'
$ws.Range("H13").WrapText = $true
$ws.Rows(13).RowHeight = 32

$ws.Range("H14").Value = 'Delicious 100% Juices combine the taste, quality and variety you love- with the convenience and Vitamin C nutrition you want. These bottles are perfect for packing in a lunch or enjoying while you are on the go.'
$ws.Range("H14").WrapText = $true
$ws.Rows(14).RowHeight = 48

$ws.Range("H15").Value = 'Super Hot Chocolate No Sugar Added Hot Cocoa Mix has 60 calories per serving and as much calcium as an 8 oz glass of milk'
$ws.Range("H15").WrapText = $true
$ws.Rows(15).RowHeight = 32

$ws.Range("H16").Value = 'Landslide Hot Chocolate No Sugar Added Hot Chocolate mix is blended with wholesome ingredients, non-fat milk from Wisconsin and premium imported cocoa'
$ws.Range("H16").WrapText = $true
$ws.Rows(16).RowHeight = 32

$ws.Range("H17").Value = 'Lightly sweetened with Stevia, an all-natural sweetener 100% Vitamin C per serving'
$ws.Range("H17").WrapText = $true

$ws.Range("H18").Value = '90% fewer calories than leading beverages; just 5 calories per 1/8 packet Sugar free Caffeine free and low sodium Kosher'
$ws.Range("H18").WrapText = $true
$ws.Rows(18).RowHeight = 32

$ws.Range("H19").Value = 'Super French Roast is in the French tradition of bold, sophisticated coffee, perfect for every morning.We start with the highest quality beans and then custom roast each bean to its peak for a full flavored taste.'
$ws.Range("H19").WrapText = $true
$ws.Rows(19).RowHeight = 48

$ws.Range("H20").Value = 'Decaf coffee without compromising the great taste you love. Super Decaf Coffee has the same aromatic flavor as Super Classic Roast Coffee, but it''s decaffeinated'
$ws.Range("H20").WrapText = $true
$ws.Rows(20).RowHeight = 32

# --- Apply wrap-text formatting (no border/fill) to the trailing blank cell H21 ---
$ws.Range("H21").WrapText = $true

# --- View state: scroll so column E is left-most, select H19 ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
$ws.Range("H19").Select()
